$wb = $excel.ActiveWorkbook

# "Reactions" sheet: rename "Flux units" header (H1) to "Flux bound units"
$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Range("H1").Value = "Flux bound units"

# "dFBA objectives" sheet: insert two new columns ("Reaction rate units",
# "Coefficient units") between "Units" (E) and "Comments" (old F)
$wsObjectives = $wb.Worksheets.Item("dFBA objectives")
$wsObjectives.Columns("F:G").Insert()
$wsObjectives.Range("F1").Value = "Reaction rate units"
$wsObjectives.Range("G1").Value = "Coefficient units"

# Make "dFBA objectives" the active sheet/tab
$wsObjectives.Activate()
